$wb = $excel.ActiveWorkbook

# --- "Front End" sheet: insert a new "Create Product" task row ---
$wsFront = $wb.Worksheets.Item("Front End")

# Insert a new row at position 8, pushing "List of Seller" / "Seller Details" down
[void]$wsFront.Rows.Item(8).Insert()

$wsFront.Cells.Item(8, 1).Value = "Create Product"
$wsFront.Cells.Item(8, 2).Value = "Seller is allowed to post new product"

# --- "Back End" sheet: extend product API description, add Admin controller row ---
$wsBack = $wb.Worksheets.Item("Back End")

# Update the product controller description to mention "create new product"
$wsBack.Cells.Item(2, 2).Value = "Implement apis for product: list product, product details,create new product"

# Add a new row for the Admin controller / approve-seller API
$wsBack.Cells.Item(5, 1).Value = "Admin controller"
$wsBack.Cells.Item(5, 2).Value = "apis for approve seller"

# --- Selection / active sheet bookkeeping to mirror the final view state ---
[void]$wsBack.Range("B3").Select()

[void]$wsFront.Activate()
[void]$wsFront.Range("C8").Select()
